$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $val) {
    $r = $ws.Range($rangeAddr)
    $r.Value = "'" + $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "47.967.98"
Set-TextValue "E2" "  +0.08%  "
Set-TextValue "D3" "2.499.53"
Set-TextValue "E3" "  -0.56%  "
Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  -0.10%  "
Set-TextValue "D5" "320.43"
Set-TextValue "E5" "  -0.82%  "
Set-TextValue "E6" "  -1.95%  "
Set-TextValue "D8" "0.999"
Set-TextValue "E8" "  -0.13%  "
Set-TextValue "E9" "  -3.46%  "
Set-TextValue "D10" "39.50"
Set-TextValue "E10" "  -3.51%  "
Set-TextValue "D11" "20.26"
Set-TextValue "E11" "  +7.89%  "
Set-TextValue "D12" "0.0811"
Set-TextValue "E12" "  -0.74%  "
Set-TextValue "E13" "  -0.35%  "
Set-TextValue "E14" "  -2.47%  "
Set-TextValue "D15" "2.889.99"
Set-TextValue "E15" "  -0.61%  "
Set-TextValue "D16" "2.487.48"
Set-TextValue "E16" "  -1.32%  "
Set-TextValue "E17" "  -2.60%  "
Set-TextValue "D18" "47.864.36"
Set-TextValue "E18" "  -0.08%  "
Set-TextValue "D19" "12.91"
Set-TextValue "E19" "  -3.19%  "
Set-TextValue "D20" "6.72"
Set-TextValue "E20" "  +0.80%  "
Set-TextValue "E21" "  -0.95%  "
Set-TextValue "D22" "2.75"
Set-TextValue "E22" "  -1.91%  "
Set-TextValue "D23" "277.94"
Set-TextValue "E23" "  +11.88%  "
Set-TextValue "D24" "71.42"
Set-TextValue "E24" "  +0.70%  "
Set-TextValue "D25" "2.54"
Set-TextValue "E25" "  -0.54%  "
Set-TextValue "E26" "  -0.08%  "
Set-TextValue "E27" "  -1.67%  "
Set-TextValue "D28" "2.10"
Set-TextValue "E28" "  -4.58%  "
Set-TextValue "E29" "  -3.29%  "
Set-TextValue "E30" "  -0.54%  "
Set-TextValue "D31" "35.00"
Set-TextValue "E31" "  -0.34%  "
Set-TextValue "D32" "49.43"
Set-TextValue "E32" "  -0.67%  "
Set-TextValue "D33" "19.47"
Set-TextValue "E33" "  -3.37%  "
Set-TextValue "E34" "  -0.17%  "
Set-TextValue "D35" "5.29"
Set-TextValue "E35" "  -1.80%  "
Set-TextValue "D36" "0.0779"
Set-TextValue "E36" "  -1.36%  "
Set-TextValue "E37" "  -2.24%  "
Set-TextValue "D38" "4.61"
Set-TextValue "E38" "  -2.18%  "
Set-TextValue "E39" "  -3.54%  "
Set-TextValue "D41" "120.66"
Set-TextValue "E41" "  +0.84%  "
Set-TextValue "E42" "  -0.30%  "
Set-TextValue "D43" "21.34"
Set-TextValue "E43" "  -5.73%  "
Set-TextValue "D44" "0.0300"
Set-TextValue "E44" "  +0.04%  "
Set-TextValue "D45" "2.005.68"
Set-TextValue "E45" "  +0.04%  "
Set-TextValue "E46" "  +1.94%  "
Set-TextValue "E47" "  -1.93%  "
Set-TextValue "E48" "  +0.02%  "
Set-TextValue "E49" "  -1.22%  "
Set-TextValue "E50" "  -1.26%  "
Set-TextValue "D51" "80.07"
Set-TextValue "E51" "  +2.46%  "
